$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the existing row 798, shifting rows 798:851 down to 800:853
$ws.Rows("798:799").Insert()

# New row 798 data
$ws.Cells.Item(798, 1).Value = 11
$ws.Cells.Item(798, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(798, 3).Value = "Bíobío"
$ws.Cells.Item(798, 4).Value = 45106
$ws.Cells.Item(798, 5).Value = 8
$ws.Cells.Item(798, 6).Value = 100112004
$ws.Cells.Item(798, 7).Value = "Cebolla"
$ws.Cells.Item(798, 8).Value = "Sin especificar"
$ws.Cells.Item(798, 9).Value = "1a (guarda)"
$ws.Cells.Item(798, 10).Value = 600
$ws.Cells.Item(798, 11).Value = 8000
$ws.Cells.Item(798, 12).Value = 8500
$ws.Cells.Item(798, 13).Value = 8250
$ws.Cells.Item(798, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(798, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(798, 16).Value = 458
$ws.Cells.Item(798, 17).Value = 18
$ws.Cells.Item(798, 18).Value = "Hortaliza"

# New row 799 data
$ws.Cells.Item(799, 1).Value = 11
$ws.Cells.Item(799, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(799, 3).Value = "Bíobío"
$ws.Cells.Item(799, 4).Value = 45106
$ws.Cells.Item(799, 5).Value = 8
$ws.Cells.Item(799, 6).Value = 100112004
$ws.Cells.Item(799, 7).Value = "Cebolla"
$ws.Cells.Item(799, 8).Value = "Sin especificar"
$ws.Cells.Item(799, 9).Value = "2a (guarda)"
$ws.Cells.Item(799, 10).Value = 300
$ws.Cells.Item(799, 11).Value = 7500
$ws.Cells.Item(799, 12).Value = 7500
$ws.Cells.Item(799, 13).Value = 7500
$ws.Cells.Item(799, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(799, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(799, 16).Value = 417
$ws.Cells.Item(799, 17).Value = 18
$ws.Cells.Item(799, 18).Value = "Hortaliza"
